$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0.3690815611690796
    3  = 0.4771425537982135
    4  = 0.3341587407122631
    5  = 0.4536773205483691
    6  = 0.3318261050981052
    7  = 0.3465738823016777
    8  = 0.4179108153570307
    9  = 0.3702067156649563
    10 = 0.4408704480605816
    11 = 0.3347437066860565
    12 = 0.4219262987243912
    13 = 0.4692907679594063
    14 = 0.4306240756002559
    15 = 0.365541752306202
    16 = 0.3930634567934926
    17 = 0.414639050402166
    18 = 0.3807190971691541
    19 = 0.5546810199847799
    20 = 0.4671580843164505
    21 = 0.331190978343616
    22 = 0.4402680554968882
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
